$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 16375.167
$ws.Range("I69").Value = 12740.333
$ws.Range("J69").Value = 20010
$ws.Range("K69").Value = 38220.999
$ws.Range("L69").Value = 60030
$ws.Range("M69").Value = -37346.999
$ws.Range("N69").Value = -61778
$ws.Range("H72").Value = 16375.167
$ws.Range("I72").Value = 12740.333
$ws.Range("J72").Value = 20010
$ws.Range("K72").Value = 114662.997
$ws.Range("L72").Value = 180090
$ws.Range("M72").Value = -110294.997
$ws.Range("N72").Value = -188826
$ws.Range("H80").Value = 1865.2593
$ws.Range("J80").Value = 3044.818
$ws.Range("L80").Value = 9134.454000000002
$ws.Range("N80").Value = -11130.454
$ws.Range("H83").Value = 1865.2593
$ws.Range("J83").Value = 3044.818
$ws.Range("L83").Value = 27403.362
$ws.Range("N83").Value = -37387.362
$ws.Range("H86").Value = 4602.4
$ws.Range("I86").Value = 1002.5
$ws.Range("K86").Value = 1002.5
$ws.Range("M86").Value = 120.5
$ws.Range("H89").Value = 4602.4
$ws.Range("I89").Value = 1002.5
$ws.Range("K89").Value = 5012.5
$ws.Range("M89").Value = 603.5
$ws.Range("H107").Value = 946.625
$ws.Range("I107").Value = 874.8570999999999
$ws.Range("K107").Value = 874.8570999999999
$ws.Range("M107").Value = 1045.1429
$ws.Range("H112").Value = 913672.0600000001
$ws.Range("J112").Value = 1004589.3
$ws.Range("L112").Value = 3013767.9
$ws.Range("N112").Value = -3015983.9
$ws.Range("H132").Value = 21743668
$ws.Range("I132").Value = 29416400
$ws.Range("K132").Value = 88249200
$ws.Range("M132").Value = -88246670
$ws.Range("H135").Value = 9538
$ws.Range("I135").Value = 6552.7144
$ws.Range("J135").Value = 14181.777
$ws.Range("K135").Value = 58974.4296
$ws.Range("L135").Value = 127635.993
$ws.Range("M135").Value = -56439.4296
$ws.Range("N135").Value = -132705.993
$ws.Range("H137").Value = 2909.7273
$ws.Range("I137").Value = 2937.739
$ws.Range("J137").Value = 2845.3
$ws.Range("K137").Value = 8813.217000000001
$ws.Range("L137").Value = 8535.900000000001
$ws.Range("M137").Value = -6263.217000000001
$ws.Range("N137").Value = -13635.9
$ws.Range("H138").Value = 1658307
$ws.Range("I138").Value = 2163.4
$ws.Range("J138").Value = 1934330.9
$ws.Range("K138").Value = 6490.200000000001
$ws.Range("L138").Value = 5802992.699999999
$ws.Range("M138").Value = -1350.200000000001
$ws.Range("N138").Value = -5813272.699999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2505.8462
$ws.Range("I122").Value = 2461.4546
$ws.Range("K122").Value = 7384.3638
$ws.Range("M122").Value = -4934.3638

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1422.3846
$ws.Range("I20").Value = 1349
$ws.Range("J20").Value = 1539.8
$ws.Range("K20").Value = 1349
$ws.Range("L20").Value = 1539.8
$ws.Range("M20").Value = -1102
$ws.Range("N20").Value = -2033.8
$ws.Range("H134").Value = 3296.4
$ws.Range("I134").Value = 3354.1052
$ws.Range("K134").Value = 10062.3156
$ws.Range("M134").Value = -7527.3156

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("H31").Value = 7382.931
$ws.Range("I31").Value = 5922.5
$ws.Range("J31").Value = 8746
$ws.Range("K31").Value = 5922.5
$ws.Range("L31").Value = 8746
$ws.Range("M31").Value = -5627.5
$ws.Range("N31").Value = -9336
$ws.Range("H34").Value = 7382.931
$ws.Range("I34").Value = 5922.5
$ws.Range("J34").Value = 8746
$ws.Range("K34").Value = 5922.5
$ws.Range("L34").Value = 8746
$ws.Range("M34").Value = -5720.5
$ws.Range("N34").Value = -9150
$ws.Range("H41").Value = 1321.2858
$ws.Range("I41").Value = 1321.2858
$ws.Range("K41").Value = 1321.2858
$ws.Range("M41").Value = -893.2858000000001
$ws.Range("H76").Value = 4340
$ws.Range("I76").Value = 4340
$ws.Range("K76").Value = 4340
$ws.Range("M76").Value = -4025
$ws.Range("H79").Value = 4340
$ws.Range("I79").Value = 4340
$ws.Range("K79").Value = 4340
$ws.Range("M79").Value = -3248
$ws.Range("H132").Value = 2161.25
$ws.Range("I132").Value = 1998.2
$ws.Range("J132").Value = 2433
$ws.Range("K132").Value = 5994.6
$ws.Range("L132").Value = 7299
$ws.Range("M132").Value = -3464.6
$ws.Range("N132").Value = -12359
$ws.Range("H134").Value = 2443.1
$ws.Range("I134").Value = 2418.8572
$ws.Range("K134").Value = 7256.571599999999
$ws.Range("M134").Value = -4721.571599999999
$ws.Range("M17").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 3113.6843
$ws.Range("J38").Value = 7048.25
$ws.Range("L38").Value = 21144.75
$ws.Range("N38").Value = -21838.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 449.66666
$ws.Range("I2").Value = 449.66666
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 449.66666
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -336.66666
$ws.Range("H80").Value = 3761.6177
$ws.Range("I80").Value = 2807.25
$ws.Range("K80").Value = 2807.25
$ws.Range("M80").Value = -1809.25
$ws.Range("H83").Value = 3761.6177
$ws.Range("I83").Value = 2807.25
$ws.Range("K83").Value = 14036.25
$ws.Range("M83").Value = -9044.25
$ws.Range("H97").Value = 1374.8572
$ws.Range("I97").Value = 876.73334
$ws.Range("K97").Value = 876.73334
$ws.Range("M97").Value = -380.73334
$ws.Range("H123").Value = 53516.8
$ws.Range("J123").Value = 53516.8
$ws.Range("L123").Value = 53516.8
$ws.Range("N123").Value = -58416.8
$ws.Range("H132").Value = 3468.4329
$ws.Range("I132").Value = 4174.2046
$ws.Range("J132").Value = 2118.261
$ws.Range("K132").Value = 12522.6138
$ws.Range("L132").Value = 6354.782999999999
$ws.Range("M132").Value = -9992.613799999999
$ws.Range("N132").Value = -11414.783
$ws.Range("H136").Value = 63427
$ws.Range("J136").Value = 63427
$ws.Range("L136").Value = 190281
$ws.Range("N136").Value = -195381
$ws.Range("N2").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9352.362999999999
$ws.Range("I40").Value = 7312.6665
$ws.Range("K40").Value = 7312.6665
$ws.Range("M40").Value = -7176.6665
$ws.Range("H70").Value = 30148
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("H73").Value = 30148
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("H82").Value = 5559.84
$ws.Range("I82").Value = 7857.625
$ws.Range("J82").Value = 1474.8889
$ws.Range("K82").Value = 7857.625
$ws.Range("L82").Value = 1474.8889
$ws.Range("M82").Value = -7496.625
$ws.Range("N82").Value = -2196.8889
$ws.Range("H85").Value = 5559.84
$ws.Range("I85").Value = 7857.625
$ws.Range("J85").Value = 1474.8889
$ws.Range("K85").Value = 7857.625
$ws.Range("L85").Value = 1474.8889
$ws.Range("M85").Value = -6609.625
$ws.Range("N85").Value = -3970.8889
$ws.Range("H93").Value = 623.16
$ws.Range("I93").Value = 428.5909
$ws.Range("K93").Value = 428.5909
$ws.Range("M93").Value = 819.4091000000001
$ws.Range("H132").Value = 7715
$ws.Range("I132").Value = 7366.6816
$ws.Range("J132").Value = 8165.7646
$ws.Range("K132").Value = 22100.0448
$ws.Range("L132").Value = 24497.2938
$ws.Range("M132").Value = -19570.0448
$ws.Range("N132").Value = -29557.2938
$ws.Range("H136").Value = 5385
$ws.Range("I136").Value = 5444.731
$ws.Range("K136").Value = 16334.193
$ws.Range("M136").Value = -13784.193
$ws.Range("N70").ClearContents()
$ws.Range("N73").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 89924.164
$ws.Range("J81").Value = 10399
$ws.Range("L81").Value = 20798
$ws.Range("N81").Value = -22920
$ws.Range("H84").Value = 89924.164
$ws.Range("J84").Value = 10399
$ws.Range("L84").Value = 103990
$ws.Range("N84").Value = -114598
$ws.Range("H107").Value = 791.8182
$ws.Range("J107").Value = 820
$ws.Range("L107").Value = 2460
$ws.Range("N107").Value = -6300
$ws.Range("H132").Value = 1215.4286
$ws.Range("I132").Value = 1162.4445
$ws.Range("J132").Value = 1533.3334
$ws.Range("K132").Value = 3487.3335
$ws.Range("L132").Value = 4600.0002
$ws.Range("M132").Value = -957.3335000000002
$ws.Range("N132").Value = -9660.0002
$ws.Range("H137").Value = 78347.5
